$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new data row (row 24) - A24 already has the date number format applied
$ws.Range("A24").Value = 45968
$ws.Range("B24").Value = 573
$ws.Range("C24").Value = 10
$ws.Range("D24").Value = 563

# Update the selection shown in the saved sheet view
$ws.Range("G21").Select()
